$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 1062.875
$ws.Range("I80").Value = 441.82608
$ws.Range("J80").Value = 2650
$ws.Range("K80").Value = 1325.47824
$ws.Range("L80").Value = 7950
$ws.Range("M80").Value = -327.4782399999999
$ws.Range("N80").Value = -9946

# Row 83
$ws.Range("H83").Value = 1062.875
$ws.Range("I83").Value = 441.82608
$ws.Range("J83").Value = 2650
$ws.Range("K83").Value = 3976.43472
$ws.Range("L83").Value = 23850
$ws.Range("M83").Value = 1015.56528
$ws.Range("N83").Value = -33834

# Row 141
$ws.Range("H141").Value = 4476.5386
$ws.Range("I141").Value = 2349.5833
$ws.Range("J141").Value = 30000
$ws.Range("K141").Value = 7048.749899999999
$ws.Range("L141").Value = 90000
$ws.Range("M141").Value = -1868.749899999999
$ws.Range("N141").Value = -100360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 4166.6665
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# Row 66
$ws.Range("H66").Value = 4166.6665
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# Row 82
$ws.Range("H82").Value = 3582
$ws.Range("I82").Value = 3582
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3582
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3221
$ws.Range("N82").ClearContents()

# Row 85
$ws.Range("H85").Value = 3582
$ws.Range("I85").Value = 3582
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3582
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2334
$ws.Range("N85").ClearContents()

# Row 132
$ws.Range("H132").Value = 2385.4583
$ws.Range("I132").Value = 1636.4615
$ws.Range("J132").Value = 3270.6365
$ws.Range("K132").Value = 4909.3845
$ws.Range("L132").Value = 9811.9095
$ws.Range("M132").Value = -2379.3845
$ws.Range("N132").Value = -14871.9095

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 25000738
$ws.Range("I86").Value = 33334018
$ws.Range("K86").Value = 33334018
$ws.Range("M86").Value = -33332895

# Row 89
$ws.Range("H89").Value = 25000738
$ws.Range("I89").Value = 33334018
$ws.Range("K89").Value = 166670090
$ws.Range("M89").Value = -166664474

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5276.4165
$ws.Range("I58").Value = 1390.5
$ws.Range("J58").Value = 7219.375
$ws.Range("K58").Value = 1390.5
$ws.Range("L58").Value = 7219.375
$ws.Range("M58").Value = -1187.5
$ws.Range("N58").Value = -7625.375

# Row 132
$ws.Range("H132").Value = 2332.6316
$ws.Range("I132").Value = 1520.091
$ws.Range("J132").Value = 3449.875
$ws.Range("K132").Value = 4560.272999999999
$ws.Range("L132").Value = 10349.625
$ws.Range("M132").Value = -2030.272999999999
$ws.Range("N132").Value = -15409.625

# Row 134
$ws.Range("H134").Value = 5940.8
$ws.Range("I134").Value = 13544.6
$ws.Range("J134").Value = 2138.9
$ws.Range("K134").Value = 40633.8
$ws.Range("L134").Value = 6416.700000000001
$ws.Range("M134").Value = -38098.8
$ws.Range("N134").Value = -11486.7

# Row 136
$ws.Range("H136").Value = 5276.4165
$ws.Range("I136").Value = 1390.5
$ws.Range("J136").Value = 7219.375
$ws.Range("K136").Value = 4171.5
$ws.Range("L136").Value = 21658.125
$ws.Range("M136").Value = -1621.5
$ws.Range("N136").Value = -26758.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 137
$ws.Range("H137").Value = 3373.2144
$ws.Range("I137").Value = 659.9091
$ws.Range("J137").Value = 5128.8823
$ws.Range("K137").Value = 1979.7273
$ws.Range("L137").Value = 15386.6469
$ws.Range("M137").Value = 3120.2727
$ws.Range("N137").Value = -25586.6469

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4523.4546
$ws.Range("I70").Value = 4545.143
$ws.Range("J70").Value = 4513.3335
$ws.Range("K70").Value = 4545.143
$ws.Range("L70").Value = 4513.3335
$ws.Range("M70").Value = -4275.143
$ws.Range("N70").Value = -5053.3335

# Row 73
$ws.Range("H73").Value = 4523.4546
$ws.Range("I73").Value = 4545.143
$ws.Range("J73").Value = 4513.3335
$ws.Range("K73").Value = 4545.143
$ws.Range("L73").Value = 4513.3335
$ws.Range("M73").Value = -3609.143
$ws.Range("N73").Value = -6385.3335

# Row 80
$ws.Range("H80").Value = 2113.3044
$ws.Range("I80").Value = 2042.8572
$ws.Range("J80").Value = 2144.125
$ws.Range("K80").Value = 2042.8572
$ws.Range("L80").Value = 2144.125
$ws.Range("M80").Value = -1044.8572
$ws.Range("N80").Value = -4140.125

# Row 83
$ws.Range("H83").Value = 2113.3044
$ws.Range("I83").Value = 2042.8572
$ws.Range("J83").Value = 2144.125
$ws.Range("K83").Value = 10214.286
$ws.Range("L83").Value = 10720.625
$ws.Range("M83").Value = -5222.286
$ws.Range("N83").Value = -20704.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 940
$ws.Range("I16").Value = 940
$ws.Range("K16").Value = 940
$ws.Range("M16").Value = -770

# Row 68
$ws.Range("H68").Value = 10254356
$ws.Range("J68").Value = 1962.8889
$ws.Range("L68").Value = 1962.8889
$ws.Range("N68").Value = -3460.8889

# Row 71
$ws.Range("H71").Value = 10254356
$ws.Range("J71").Value = 1962.8889
$ws.Range("L71").Value = 9814.4445
$ws.Range("N71").Value = -17302.4445

# Row 82
$ws.Range("H82").Value = 5348907.5
$ws.Range("I82").Value = 12988000
$ws.Range("J82").Value = 1542.9
$ws.Range("K82").Value = 12988000
$ws.Range("L82").Value = 1542.9
$ws.Range("M82").Value = -12987639
$ws.Range("N82").Value = -2264.9

# Row 85
$ws.Range("H85").Value = 5348907.5
$ws.Range("I85").Value = 12988000
$ws.Range("J85").Value = 1542.9
$ws.Range("K85").Value = 12988000
$ws.Range("L85").Value = 1542.9
$ws.Range("M85").Value = -12986752
$ws.Range("N85").Value = -4038.9

# Row 100
$ws.Range("H100").Value = 1725.5333
$ws.Range("I100").Value = 1613.8334
$ws.Range("J100").Value = 1800
$ws.Range("K100").Value = 1613.8334
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -1072.8334
$ws.Range("N100").Value = -2882

# Row 132
$ws.Range("H132").Value = 10210535
$ws.Range("I132").Value = 27791278
$ws.Range("J132").Value = 2362.2258
$ws.Range("K132").Value = 83373834
$ws.Range("L132").Value = 7086.6774
$ws.Range("M132").Value = -83371304
$ws.Range("N132").Value = -12146.6774

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1866.9333
$ws.Range("J136").Value = 2190.9092
$ws.Range("L136").Value = 6572.7276
$ws.Range("N136").Value = -11672.7276
